# Auto-generated: apply price/volume updates to Sheet1 (rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.104.91'
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("D3").Value = '1.667.34'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5114'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.23%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2637'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.86%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06429'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07413'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.85%  '
$ws.Range("D12").Value = '1.667.37'
$ws.Range("E12").Value = '  -2.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.508'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5817'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008574'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.44'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("D17").Value = '26.157.13'
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.927'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.214'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.005'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.639'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1197'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.46%  '
$ws.Range("E27").Value = '  +1.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06464'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.305'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.319'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.521'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.515'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.637'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.018'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6083'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.57%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.668'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.205'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.74%  '
$ws.Range("E39").Value = '  +2.13%  '
$ws.Range("D40").Value = '1.078.39'
$ws.Range("E40").Value = '  +0.95%  '
$ws.Range("E41").Value = '  +1.63%  '
$ws.Range("E42").Value = '  +0.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.77'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.85%  '
$ws.Range("D44").Value = '1.816.42'
$ws.Range("E44").Value = '  -0.99%  '
$ws.Range("E45").Value = '  +9.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.064'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05207'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("E50").Value = '  -0.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.958'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.49%  '
